$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Bl"
$ws.Range("G1").Value = "Operating Freq"

$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

$ws.Range("J11").Select()
